$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 649.5
$ws.Range("I31").Value = 649.5
$ws.Range("K31").Value = 1948.5
$ws.Range("M31").Value = -1718.5

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3822.6365
$ws.Range("I32").Value = 2871.6
$ws.Range("K32").Value = 2871.6
$ws.Range("M32").Value = -2584.6

$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -3812

$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -5808

$ws.Range("H122").Value = 4142.933
$ws.Range("I122").Value = 4301.6924
$ws.Range("J122").Value = 3111
$ws.Range("K122").Value = 12905.0772
$ws.Range("L122").Value = 9333
$ws.Range("M122").Value = -10455.0772
$ws.Range("N122").Value = -14233

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4400
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4600
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4600
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6846

$ws.Range("H89").Value = 4400
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4600
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 23000
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -34232

$ws.Range("H92").Value = 139499.75
$ws.Range("J92").Value = 139499.75
$ws.Range("L92").Value = 139499.75
$ws.Range("N92").Value = -144491.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 21988
$ws.Range("J28").Value = 21988
$ws.Range("L28").Value = 21988
$ws.Range("N28").Value = -22478

$ws.Range("H32").Value = 3732.5
$ws.Range("I32").Value = 2478.8
$ws.Range("K32").Value = 2478.8
$ws.Range("M32").Value = -2162.8

$ws.Range("H38").Value = 31800
$ws.Range("J38").Value = 50000
$ws.Range("L38").Value = 50000
$ws.Range("N38").Value = -50754

$ws.Range("H46").Value = 31800
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50422

$ws.Range("H51").Value = 23363.334
$ws.Range("I51").Value = 23363.334
$ws.Range("K51").Value = 23363.334
$ws.Range("M51").Value = -22627.334

$ws.Range("H61").Value = 23363.334
$ws.Range("I61").Value = 23363.334
$ws.Range("K61").Value = 23363.334
$ws.Range("M61").Value = -23015.334

$ws.Range("H132").Value = 1750
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H141").Value = 639994.5600000001
$ws.Range("J141").Value = 639994.5600000001
$ws.Range("L141").Value = 639994.5600000001
$ws.Range("N141").Value = -650354.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 719.3
$ws.Range("I6").Value = 38.2
$ws.Range("J6").Value = 1400.4
$ws.Range("K6").Value = 114.6
$ws.Range("L6").Value = 4201.200000000001
$ws.Range("M6").Value = -1.600000000000009
$ws.Range("N6").Value = -4427.200000000001

$ws.Range("H17").Value = 1069.8
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1212.25
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 3636.75
$ws.Range("N17").Value = -3974.75
$ws.Range("M17").Value = -1331

$ws.Range("H34").Value = 1620
$ws.Range("I34").Value = 149
$ws.Range("J34").Value = 2600.6667
$ws.Range("K34").Value = 447
$ws.Range("L34").Value = 7802.000100000001
$ws.Range("M34").Value = -363
$ws.Range("N34").Value = -7970.000100000001

$ws.Range("H39").Value = 1003
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 4021.7424
$ws.Range("J55").Value = 4073.6155
$ws.Range("L55").Value = 12220.8465
$ws.Range("N55").Value = -12574.8465

$ws.Range("I109").Value = 450.66666
$ws.Range("J109").Value = 500
$ws.Range("K109").Value = 1351.99998
$ws.Range("L109").Value = 1500
$ws.Range("M109").Value = -311.9999800000001
$ws.Range("N109").Value = -3580

$ws.Range("H120").Value = 8555.75
$ws.Range("I120").Value = 1407.6666
$ws.Range("J120").Value = 30000
$ws.Range("K120").Value = 4222.9998
$ws.Range("L120").Value = 90000
$ws.Range("M120").Value = 615.0002000000004
$ws.Range("N120").Value = -99676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3379.7
$ws.Range("I97").Value = 3541.6
$ws.Range("J97").Value = 3217.8
$ws.Range("K97").Value = 3541.6
$ws.Range("L97").Value = 3217.8
$ws.Range("M97").Value = -3045.6
$ws.Range("N97").Value = -4209.8

$ws.Range("H98").Value = 12918.6
$ws.Range("J98").Value = 12918.6
$ws.Range("L98").Value = 12918.6
$ws.Range("N98").Value = -18908.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 38499.5
$ws.Range("J68").Value = 38499.5
$ws.Range("L68").Value = 38499.5
$ws.Range("N68").Value = -40121.5

$ws.Range("H71").Value = 38499.5
$ws.Range("J71").Value = 38499.5
$ws.Range("L71").Value = 115498.5
$ws.Range("N71").Value = -123610.5

$ws.Range("H112").Value = 32693.5
$ws.Range("J112").Value = 32693.5
$ws.Range("L112").Value = 32693.5
$ws.Range("N112").Value = -35647.5

$ws.Range("H132").Value = 1828.4667
$ws.Range("I132").Value = 1754.3334
$ws.Range("J132").Value = 2125
$ws.Range("K132").Value = 5263.0002
$ws.Range("L132").Value = 6375
$ws.Range("M132").Value = -2733.0002
$ws.Range("N132").Value = -11435

